$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the value for the removed URL rows (A3:A6), keeping their style.
$ws.Range("A3:A6").ClearContents()

# Rebuild the hyperlinks collection so only A2's hyperlink survives
# (per-hyperlink .Delete() isn't effective, so drop them all and re-add A2's).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://www.wineenthusiast.com/") | Out-Null
$ws.Range("A2").Style = "Hyperlink"

# Update the current selection to A6 (matches the saved sheetView selection)
$ws.Range("A6").Select()
